$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.806.59'
$ws.Range('E2').Value = '  -1.63%  '

$ws.Range('D3').Value = '1.872.64'
$ws.Range('E3').Value = '  -2.41%  '

$ws.Range('E4').Value = '  -0.70%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.93'
$ws.Range('E5').Value = '  -2.22%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.689'
$ws.Range('E6').Value = '  -6.38%  '

$ws.Range('E7').Value = '  -0.81%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.08'
$ws.Range('E8').Value = '  +3.37%  '

$ws.Range('E9').Value = '  -3.03%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '50.92'
$ws.Range('E10').Value = '  -3.31%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0736'
$ws.Range('E11').Value = '  -0.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0969'
$ws.Range('E12').Value = '  -2.99%  '

$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '12.86'
$ws.Range('E13').Value = '  +1.87%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.144.19'
$ws.Range('E14').Value = '  -2.19%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.714'
$ws.Range('E15').Value = '  -0.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.89'
$ws.Range('E16').Value = '  -0.32%  '

$ws.Range('D17').Value = '1.885.13'
$ws.Range('E17').Value = '  -1.87%  '

$ws.Range('D18').Value = '34.784.21'
$ws.Range('E18').Value = '  -1.64%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.75'
$ws.Range('E19').Value = '  -0.56%  '

$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('E20').Value = '  -1.82%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.77'
$ws.Range('E21').Value = '  +0.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.72'
$ws.Range('E22').Value = '  -2.88%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.92'
$ws.Range('E23').Value = '  -3.07%  '

$ws.Range('E24').Value = '  -0.81%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.44'
$ws.Range('E25').Value = '  +4.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  -3.04%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.24'
$ws.Range('E27').Value = '  -1.53%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.37'
$ws.Range('E28').Value = '  -4.18%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.20'
$ws.Range('E29').Value = '  -3.32%  '

$ws.Range('E30').Value = '  -6.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.70'
$ws.Range('E32').Value = '  +3.37%  '

$ws.Range('E33').Value = '  -1.95%  '

$ws.Range('E34').Value = '  -0.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.15'
$ws.Range('E35').Value = '  -2.77%  '

$ws.Range('E36').Value = '  -0.75%  '

$ws.Range('E37').Value = '  -9.11%  '

$ws.Range('E38').Value = '  -2.61%  '

$ws.Range('E39').Value = '  -22.91%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '97.46'
$ws.Range('E40').Value = '  -1.62%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.86'
$ws.Range('E41').Value = '  -3.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0658'
$ws.Range('E42').Value = '  +1.19%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0210'
$ws.Range('E43').Value = '  -0.14%  '

$ws.Range('E44').Value = '  -5.30%  '

$ws.Range('D45').Value = '1.282.46'
$ws.Range('E45').Value = '  -4.91%  '

$ws.Range('E46').Value = '  -5.96%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0792'
$ws.Range('E47').Value = '  +8.67%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('E48').Value = '  -1.28%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.73'
$ws.Range('E49').Value = '  -2.26%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.05'
$ws.Range('E50').Value = '  +3.46%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.44'
$ws.Range('E51').Value = '  -4.16%  '
